# Solve Leetcode - 355. Design Twitter - Linked Lists and Heaps
# Adds a new row (16) to the "Neetcode 150" sheet for the "355. Design Twitter"
# problem, filed under the Heap/PQ section (notes mention Linked Lists too).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New row 16 data -------------------------------------------------
# Column A normally holds the CATEGORY, but for this entry the note
# "Heap/PQ / Linked Lists" (two lines) was typed there instead.
$ws.Cells.Item(16, 1).Value = "Heap/PQ`nLinked Lists"
$ws.Cells.Item(16, 1).WrapText = $true

# Column B = DIFFICULTY
$ws.Cells.Item(16, 2).Value = "Medium"

# Column C = NAME, with a hyperlink to the Leetcode problem page.
$link = $ws.Hyperlinks.Add(
    $ws.Cells.Item(16, 3),
    "https://leetcode.com/problems/design-twitter/",
    [Type]::Missing,
    [Type]::Missing,
    "https://leetcode.com/problems/design-twitter/"
)
$ws.Cells.Item(16, 3).Value = "355. Design Twitter"

# Match the row height used by the other two-line wrapped rows (e.g. row 2).
$ws.Rows.Item(16).RowHeight = 28.8

# --- View state --------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D16").Select()
